$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 189, pushing the old rows
# 189-241 down to 193-245 (their content/styles are preserved by Excel's
# native row-insert shift).
$ws.Rows("189:192").Insert()

# Row 189 (new): Hass / Especial
$ws.Range("A189").Value = 1
$ws.Range("B189").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C189").Value = "Arica y Parinacota"
$ws.Range("D189").Value = 45135
$ws.Range("E189").Value = 15
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100106
$ws.Range("H189").Value = "Oleaginosos"
$ws.Range("I189").Value = 100106002
$ws.Range("J189").Value = "Palta"
$ws.Range("K189").Value = "Hass"
$ws.Range("L189").Value = "Especial"
$ws.Range("M189").Value = 104
$ws.Range("N189").Value = 29000
$ws.Range("O189").Value = 30000
$ws.Range("P189").Value = 29500
$ws.Range("Q189").Value = "$/bandeja 10 kilos"
$ws.Range("R189").Value = "Perú"
$ws.Range("S189").Value = 2950
$ws.Range("T189").Value = 10

# Row 190 (new): Hass / Primera
$ws.Range("A190").Value = 1
$ws.Range("B190").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C190").Value = "Arica y Parinacota"
$ws.Range("D190").Value = 45135
$ws.Range("E190").Value = 15
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100106
$ws.Range("H190").Value = "Oleaginosos"
$ws.Range("I190").Value = 100106002
$ws.Range("J190").Value = "Palta"
$ws.Range("K190").Value = "Hass"
$ws.Range("L190").Value = "Primera"
$ws.Range("M190").Value = 208
$ws.Range("N190").Value = 27000
$ws.Range("O190").Value = 28000
$ws.Range("P190").Value = 27500
$ws.Range("Q190").Value = "$/bandeja 10 kilos"
$ws.Range("R190").Value = "Perú"
$ws.Range("S190").Value = 2750
$ws.Range("T190").Value = 10

# Row 191 (new): Hass / Segunda
$ws.Range("A191").Value = 1
$ws.Range("B191").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C191").Value = "Arica y Parinacota"
$ws.Range("D191").Value = 45135
$ws.Range("E191").Value = 15
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100106
$ws.Range("H191").Value = "Oleaginosos"
$ws.Range("I191").Value = 100106002
$ws.Range("J191").Value = "Palta"
$ws.Range("K191").Value = "Hass"
$ws.Range("L191").Value = "Segunda"
$ws.Range("M191").Value = 208
$ws.Range("N191").Value = 25000
$ws.Range("O191").Value = 26000
$ws.Range("P191").Value = 25500
$ws.Range("Q191").Value = "$/bandeja 10 kilos"
$ws.Range("R191").Value = "Perú"
$ws.Range("S191").Value = 2550
$ws.Range("T191").Value = 10

# Row 192 (new): Hass / Tercera
$ws.Range("A192").Value = 1
$ws.Range("B192").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C192").Value = "Arica y Parinacota"
$ws.Range("D192").Value = 45135
$ws.Range("E192").Value = 15
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100106
$ws.Range("H192").Value = "Oleaginosos"
$ws.Range("I192").Value = 100106002
$ws.Range("J192").Value = "Palta"
$ws.Range("K192").Value = "Hass"
$ws.Range("L192").Value = "Tercera"
$ws.Range("M192").Value = 104
$ws.Range("N192").Value = 24000
$ws.Range("O192").Value = 24000
$ws.Range("P192").Value = 24000
$ws.Range("Q192").Value = "$/bandeja 10 kilos"
$ws.Range("R192").Value = "Perú"
$ws.Range("S192").Value = 2400
$ws.Range("T192").Value = 10

Write-Host "Applied insert of 4 new rows at 189-192."
